# Update cryptos list figures (prices / 1h volume %) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.573.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.628.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +2.55%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0621"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.36%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.857.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.640.64"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.558.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.89%  "
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.32"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  +4.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("E28").Value = "  +3.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0502"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.97%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.217.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.21%  "
$ws.Range("E37").Value = "  +5.07%  "
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.498"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.23%  "
$ws.Range("E41").Value = "  -1.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.793"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.33"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.763.99"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.76%  "
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.410"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
